# Remove the "「お問い合わせください」" post row (row 384). All rows below it
# (385-496) shift up by one, and the used range shrinks from A1:C496 to A1:C495.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).EntireRow.Delete()
